{"js": "// Helper: find a unique body-text match and replace it with new text.\nasync function replaceOnce(body, searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + searchText);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1) \"Existing data\" paragraph: the pilot-study / main-experiment data collection\n//    sentence is reworded (\"collected\" -> \"collected and analyzed ... collected\").\nawait replaceOnce(\n  body,\n  \"we have collected data from a small pilot study and four data sets from the main experiment, due to feedback-induced changes in the design.\",\n  \"we have collected and analyzed data from a small pilot study and collected data from the main experiment, due to feedback-induced changes in the design.\"\n);\n\n// 2) The following parenthetical now refers back to \"the latter\" (main experiment)\n//    instead of repeating \"These four data sets\".\nawait replaceOnce(\n  body,\n  \" (These four data sets will be included in the final analysis,\",\n  \" (Four data sets of the latter will be included in the final analysis,\"\n);\n\n// 3) Stopping rule: \"noon\" -> \"TIME\" placeholder, \"12th\" -> \"17th\" of August.\nawait replaceOnce(body, \" noon of the 12\", \" TIME of the 17\");\n\n// 4) \"...which is 12 days...\" -> \"...which is 17 days...\"\nawait replaceOnce(body, \" of August, which is 12\", \" of August, which is 17\");\n\n// 5) Append a new explanatory sentence at the end of the Stopping rule paragraph.\nawait replaceOnce(\n  body,\n  \"s after sending out the invitations, due to the project deadline.\",\n  \"s after sending out the invitations, due to the project deadline. (This was originally planned for 12 days, but due to technical problems and feedback, we extended the data collection period.)\"\n);\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"Existing data\" paragraph: the pilot-study / main-experiment data collection\n#    sentence is reworded (\"collected\" -> \"collected and analyzed ... collected\").\n$rng = $d.Content\n$rng.Find.Execute(\"we have collected data from a small pilot study and four data sets from the main experiment, due to feedback-induced changes in the design.\") | Out-Null\n$rng.Text = \"we have collected and analyzed data from a small pilot study and collected data from the main experiment, due to feedback-induced changes in the design.\"\n\n# 2) The following parenthetical now refers back to \"the latter\" (main experiment)\n#    instead of repeating \"These four data sets\".\n$rng = $d.Content\n$rng.Find.Execute(\" (These four data sets will be included in the final analysis,\") | Out-Null\n$rng.Text = \" (Four data sets of the latter will be included in the final analysis,\"\n\n# 3) Stopping rule: \"noon\" -> \"TIME\" placeholder, \"12th\" -> \"17th\" of August.\n$rng = $d.Content\n$rng.Find.Execute(\" noon of the 12\") | Out-Null\n$rng.Text = \" TIME of the 17\"\n\n# 4) \"...which is 12 days...\" -> \"...which is 17 days...\"\n$rng = $d.Content\n$rng.Find.Execute(\" of August, which is 12\") | Out-Null\n$rng.Text = \" of August, which is 17\"\n\n# 5) Append a new explanatory sentence at the end of the Stopping rule paragraph.\n$rng = $d.Content\n$rng.Find.Execute(\"s after sending out the invitations, due to the project deadline.\") | Out-Null\n$rng.Text = \"s after sending out the invitations, due to the project deadline. (This was originally planned for 12 days, but due to technical problems and feedback, we extended the data collection period.)\"\n"}
